$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the ZIP value for row 2 (John Smith) and the City value for row 4 (Brad Pitt)
$ws.Range("F2").ClearContents()
$ws.Range("D4").ClearContents()

# Update the active selection to F2
$ws.Range("F2").Select()
